$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.157.27"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.51%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.600.25"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.64%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.79"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.70%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +3.00%  "
$ws.Range("E9").Value = "  +3.79%  "
$ws.Range("E10").Value = "  +3.72%  "
$ws.Range("E11").Value = "  +0.20%  "
$ws.Range("E12").Value = "  +0.84%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.38"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.14%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.065.03"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.83%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "63.012.07"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.40%  "
$ws.Range("E16").Value = "  +4.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.599.57"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.64%  "
$ws.Range("E18").Value = "  +1.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "344.36"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.53%  "
$ws.Range("E20").Value = "  +2.86%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.71"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.46%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "67.24"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.723.30"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.82%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.169"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.60"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.996"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.42"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.90%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.87"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +9.01%  "
$ws.Range("E31").Value = "  +0.79%  "
$ws.Range("E32").Value = "  +5.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "471.27"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +19.31%  "
$ws.Range("E34").Value = "  +2.93%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "176.49"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.22%  "
$ws.Range("E36").Value = "  +6.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.405"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("E39").Value = "  +1.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.59"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.99%  "
$ws.Range("E42").Value = "  -1.17%  "
$ws.Range("E43").Value = "  +6.12%  "
$ws.Range("E44").Value = "  +2.69%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.642"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.77%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "21.18"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.28%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0549"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.32%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0973"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.35%  "
$ws.Range("E49").Value = "  +1.38%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.67"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.78%  "
$ws.Range("E51").Value = "  +3.56%  "
